$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in column D and C as per the recalculated result data
$ws.Range("D3").Value = -5.926399999999998
$ws.Range("D4").Value = -7.824100000000004
$ws.Range("D7").Value = -7.918999999999994
$ws.Range("D8").Value = -8.662499999999994
$ws.Range("C11").Value = -13.3969
$ws.Range("C12").Value = -14.56430000000003
$ws.Range("D12").Value = -8.184400000000002
$ws.Range("D14").Value = -8.581599999999995
$ws.Range("C15").Value = -11.9256
$ws.Range("D22").Value = -8.335699999999997
